$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---

# Cells removed entirely in the update
$ws.Range("D2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("L2").ClearContents()

# Updated values
$ws.Range("K2").Value = -29.6
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 63.4
$ws.Range("V2").Value = 0.4446002805049089
$ws.Range("W2").Value = -0.5304659498207885
$ws.Range("X2").Value = 0.08085818873074019
$ws.Range("Y2").Value = -0.6113241385515287
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = -0.4142857142857144
$ws.Range("AB2").Value = 0.07025180545980589
$ws.Range("AC2").Value = -0.4845375197455203
$ws.Range("AD2").Value = 52.2
$ws.Range("AF2").Value = 52.2
$ws.Range("AG2").Value = -11.2
$ws.Range("AH2").Value = 0.2679671457905544
$ws.Range("AI2").Value = 0.6761658031088082
$ws.Range("AJ2").Value = -0.08523592085235918
$ws.Range("AK2").Value = -0.8115942028985501
$ws.Range("AM2").Value = -1.02
$ws.Range("AN2").Value = -9.648798521256932
$ws.Range("AP2").Value = 2.070240295748613
$ws.Range("AQ2").Value = 5.401960784313725

# --- Row 3 ---

# Cells removed entirely in the update
$ws.Range("D3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("L3").ClearContents()

# Updated values
$ws.Range("K3").Value = -29.6
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 63.4
$ws.Range("V3").Value = 0.4446002805049089
$ws.Range("W3").Value = -0.5304659498207885
$ws.Range("X3").Value = 0.08085818873074019
$ws.Range("Y3").Value = -0.6113241385515287
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = -0.4142857142857144
$ws.Range("AB3").Value = 0.07025180545980589
$ws.Range("AC3").Value = -0.4845375197455203
$ws.Range("AD3").Value = 52.2
$ws.Range("AF3").Value = 52.2
$ws.Range("AG3").Value = -11.2
$ws.Range("AH3").Value = 0.2679671457905544
$ws.Range("AI3").Value = 0.6761658031088082
$ws.Range("AJ3").Value = -0.08523592085235918
$ws.Range("AK3").Value = -0.8115942028985501
$ws.Range("AM3").Value = -1.02
$ws.Range("AN3").Value = -9.648798521256932
$ws.Range("AP3").Value = 2.070240295748613
$ws.Range("AQ3").Value = 5.401960784313725
